# Update timing parameters in row 4 of Sheet1 (Aristoteles row)
# per commit "Parametros para menos de 18 curvas"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H4").Value = 1140
$ws.Range("I4").Value = 1023
$ws.Range("J4").Value = 1099
$ws.Range("Q4").Value = 705
